$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "309.24"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-3.74%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "50.35"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "3.26%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.198"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.19%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07767"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-4.33%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.496"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-2.08%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.340"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "10.98%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.564"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-5.05%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-6.35%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1975"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.42%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.04792"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "3.62%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09432"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.69%"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.55%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001268"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-4.55%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005794"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-2.41%"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "2,016.37%"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.20%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.437"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.28%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3477"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2.14%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.985"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.08%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1367"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.13%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04166"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.00%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001269"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-2.97%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.003937"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-7.48%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001349"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.04%"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-4.21%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06018"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "4.22%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01099"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "74.43%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007978"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "3.63%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1423"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-1.50%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.008387"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "8.98%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008340"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "2.94%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3367"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "5.44%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00007206"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "2.61%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.02%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-3.31%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002619"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-34.52%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.02%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.02%"
